$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update the "Model" column (C) values — Level1/Level2Phero counts bumped,
# and the repeated Level3Phero rows raised from 20 to 50.
$ws.Range("C2").Value = 10
$ws.Range("C3").Value = 20
$ws.Range("C4:C10").Value = 50

# Selection moved from F7 to L10 (last active cell in the sheet view).
$ws.Range("L10").Select()
